$d = $word.ActiveDocument

# 1) Replace the first citation (spans three runs, middle one a red space)
#    with the updated decision number/date/ADA text, merged into a single run.
$d.Content.Find.Execute(
    "Την με αριθ. Φ.350.2/1/32958/Ε3/27-2-2018  (ΑΔΑ:6Π414653ΠΣ-7ΕΝ) Υπουργική Απόφαση με θέμα: «Τοποθέτηση Περιφερειακών Διευθυντών Εκπαίδευσης »",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Τη με αριθ. Φ.351.1/11/48020/Ε3/28-3-2019 (ΑΔΑ: ΩΩΤΗ4653ΠΣ-ΒΔ3) Υπουργική Απόφαση με θέμα: «Τοποθέτηση Περιφερειακών Διευθυντών Εκπαίδευσης»",
    2)

# 2) Fix grammatical form "Την" -> "Τη" before the second citation.
$d.Content.Find.Execute(
    "Την με αριθ. Φ.353.1/324/105657/Δ1/8-10-2002",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Τη με αριθ. Φ.353.1/324/105657/Δ1/8-10-2002",
    2)
